$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain plain text, matching the source export
# which always stores prices as inline strings (not numbers).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.137.99'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.29%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.612.62'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.85%  '

$ws.Range("E4").Value = '  -0.41%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.00'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.03%  '

$ws.Range("E6").Value = '  -0.42%  '

$ws.Range("E7").Value = '  +1.18%  '

$ws.Range("E8").Value = '  +1.67%  '

$ws.Range("E9").Value = '  +1.86%  '

$ws.Range("E10").Value = '  +3.30%  '

$ws.Range("E11").Value = '  +1.51%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.836.34'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.86%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.612.16'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.91%  '

$ws.Range("E14").Value = '  +0.47%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.512'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.64%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.144.61'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.33%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.94'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.06%  '

$ws.Range("E18").Value = '  +2.02%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '198.97'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.25%  '

$ws.Range("E21").Value = '  +2.36%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.49'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.31%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.04'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.79%  '

$ws.Range("E24").Value = '  +3.31%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '142.68'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.09%  '

$ws.Range("E26").Value = '  +2.45%  '

$ws.Range("E27").Value = '  -0.48%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.22'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.43%  '

$ws.Range("E29").Value = '  +0.11%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.18'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.63%  '

$ws.Range("E31").Value = '  +2.85%  '

$ws.Range("E32").Value = '  +2.38%  '

$ws.Range("E33").Value = '  +1.75%  '

$ws.Range("E34").Value = '  +3.71%  '

$ws.Range("E35").Value = '  -1.90%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.108.74'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.11%  '

$ws.Range("E37").Value = '  +1.64%  '

$ws.Range("E38").Value = '  -0.48%  '

$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.507'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.08%  '

$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.34'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.10%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.792'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.39%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.797'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.93%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.749.07'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.92%  '

$ws.Range("E44").Value = '  +1.47%  '

$ws.Range("E45").Value = '  -2.56%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₆0107'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +9.38%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.56'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +9.86%  '

$ws.Range("E48").Value = '  +1.73%  '

$ws.Range("E49").Value = '  +0.08%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.407'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.75%  '

$ws.Range("E51").Value = '  -0.44%  '

